$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.469.40"
$ws.Range("E2").Value = "  +3.35%  "
$ws.Range("D3").Value = "1.915.20"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.61"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.93"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.36"
$ws.Range("E9").Value = "  +8.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.365"
$ws.Range("E10").Value = "  +2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0764"
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0994"
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.42"
$ws.Range("E13").Value = "  +8.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.799"
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("D15").Value = "2.194.72"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("E16").Value = "  +4.68%  "
$ws.Range("D17").Value = "1.914.49"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "36.376.68"
$ws.Range("E18").Value = "  +2.67%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "251.95"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.16"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  +4.43%  "
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.87"
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.74"
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.84"
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("E31").Value = "  +5.96%  "
$ws.Range("E32").Value = "  +4.21%  "
$ws.Range("E33").Value = "  +6.77%  "
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0851"
$ws.Range("E36").Value = "  +22.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("E37").Value = "  -15.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.860"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("E39").Value = "  +2.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.80"
$ws.Range("E40").Value = "  +9.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0229"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.87"
$ws.Range("E42").Value = "  +28.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.07"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "1.341.19"
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.78"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.43"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("D51").Value = "2.093.00"
$ws.Range("E51").Value = "  +1.19%  "
